$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.774.32'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '1.625.05'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''214.94'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').Value = '''0.5105'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D8').Value = '''0.2562'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '''0.06321'
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('D10').Value = '''19.35'
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('D11').Value = '''0.07779'
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.631.69'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.221'
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').Value = '1.846.85'
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').Value = '''0.5532'
$ws.Range('E15').Value = '  +2.24%  '
$ws.Range('D16').Value = '''63.43'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').Value = '0.0₅7502'
$ws.Range('E17').Value = '  -2.23%  '
$ws.Range('D18').Value = '25.781.36'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').Value = '''4.407'
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').Value = '''193.73'
$ws.Range('E21').Value = '  -2.70%  '
$ws.Range('D22').Value = '''9.745'
$ws.Range('E22').Value = '  -1.39%  '
$ws.Range('D23').Value = '''5.995'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('D24').Value = '''1.002'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = '''1.868'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').Value = '''141.52'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = '''0.1245'
$ws.Range('E27').Value = '  +4.01%  '
$ws.Range('D28').Value = '''15.48'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('D29').Value = '''6.701'
$ws.Range('E29').Value = '  -1.63%  '
$ws.Range('D30').Value = '''1.236'
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('D31').Value = '''0.04841'
$ws.Range('E31').Value = '  -1.26%  '
$ws.Range('D32').Value = '''3.235'
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('D33').Value = '''3.162'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('D35').Value = '''2.365'
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').Value = '''0.8911'
$ws.Range('E36').Value = '  -1.95%  '
$ws.Range('D37').Value = '''2.538'
$ws.Range('E37').Value = '  -1.75%  '
$ws.Range('E38').Value = '  +0.94%  '
$ws.Range('D39').Value = '1.114.14'
$ws.Range('E39').Value = '  -2.23%  '
$ws.Range('D40').Value = '''0.01544'
$ws.Range('E40').Value = '  -1.24%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').Value = '''5.512'
$ws.Range('E42').Value = '  +1.79%  '
$ws.Range('D43').Value = '''0.7948'
$ws.Range('E43').Value = '  -1.66%  '
$ws.Range('D44').Value = '''96.99'
$ws.Range('E44').Value = '  -1.95%  '
$ws.Range('D45').Value = '1.771.79'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('E46').Value = '  -8.25%  '
$ws.Range('D47').Value = '''0.4423'
$ws.Range('E47').Value = '  -2.36%  '
$ws.Range('D48').Value = '''0.9979'
$ws.Range('E48').Value = '  -0.74%  '
$ws.Range('D49').Value = '''54.51'
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('D50').Value = '''0.05127'
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('D51').Value = '''7.530'
$ws.Range('E51').Value = '  +2.76%  '
